$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Update the existing value for period 1 (was 5558000000) to reflect the
# new temporal resolution, then extend the series through period 12,
# copying period 1's formatting down the new rows.
$ws.Cells.Item(3, 2).Value = 463166667
$ws.Cells.Item(3, 2).Copy()

for ($row = 4; $row -le 14; $row++) {
    $period = $row - 2
    $ws.Cells.Item($row, 1).Value = $period
    $ws.Cells.Item($row, 2).Value = 463166667
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
}

# Column B was widened to fit the new values.
$ws.Columns.Item(2).ColumnWidth = 9.17

# The Demand sheet is now the active tab / selected cell.
$ws.Activate()
$ws.Range("H13").Select()
